$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.891504666666667
$ws.Range("H2").Value = 8.674514
$ws.Range("I2").Value = 0.1213590456377548
$ws.Range("J2").Value = 0.1213590456377548
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 17.16653
$ws.Range("N2").Value = 51.49959
$ws.Range("O2").Value = 0.0560345397128279
$ws.Range("P2").Value = 0.0560345397128279
$ws.Range("Q2").Value = 49.63710160547333
$ws.Range("R2").Value = 446.73391444926
$ws.Range("S2").Value = 0.006800298262299663
$ws.Range("T2").Value = 0.006800298262299663

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.891504666666667
$ws.Range("H3").Value = 8.674514
$ws.Range("I3").Value = 0.1213590456377548
$ws.Range("J3").Value = 0.1213590456377548
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 256.4443053333333
$ws.Range("N3").Value = 769.332916
$ws.Range("O3").Value = 0.8370788162388805
$ws.Range("P3").Value = 0.8370788162388805
$ws.Range("Q3").Value = 741.5099056114249
$ws.Range("R3").Value = 6673.589150502824
$ws.Range("S3").Value = 0.101587086262332
$ws.Range("T3").Value = 0.101587086262332

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.891504666666667
$ws.Range("H4").Value = 8.674514
$ws.Range("I4").Value = 0.1213590456377548
$ws.Range("J4").Value = 0.1213590456377548
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 32.74538866666666
$ws.Range("N4").Value = 98.236166
$ws.Range("O4").Value = 0.1068866440482915
$ws.Range("P4").Value = 0.1068866440482915
$ws.Range("Q4").Value = 94.68344414148044
$ws.Range("R4").Value = 852.150997273324
$ws.Range("S4").Value = 0.01297166111312306
$ws.Range("T4").Value = 0.01297166111312306

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 12.04042966666667
$ws.Range("H5").Value = 36.121289
$ws.Range("I5").Value = 0.505347637947847
$ws.Range("J5").Value = 0.505347637947847
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.16653
$ws.Range("N5").Value = 51.49959
$ws.Range("O5").Value = 0.0560345397128279
$ws.Range("P5").Value = 0.0560345397128279
$ws.Range("Q5").Value = 206.6923970857233
$ws.Range("R5").Value = 1860.23157377151
$ws.Range("S5").Value = 0.0283169222873724
$ws.Range("T5").Value = 0.0283169222873724

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.04042966666667
$ws.Range("H6").Value = 36.121289
$ws.Range("I6").Value = 0.505347637947847
$ws.Range("J6").Value = 0.505347637947847
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 256.4443053333333
$ws.Range("N6").Value = 769.332916
$ws.Range("O6").Value = 0.8370788162388805
$ws.Range("P6").Value = 0.8370788162388805
$ws.Range("Q6").Value = 3087.699621783192
$ws.Range("R6").Value = 27789.29659604873
$ws.Range("S6").Value = 0.4230158025624981
$ws.Range("T6").Value = 0.4230158025624981

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.04042966666667
$ws.Range("H7").Value = 36.121289
$ws.Range("I7").Value = 0.505347637947847
$ws.Range("J7").Value = 0.505347637947847
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 32.74538866666666
$ws.Range("N7").Value = 98.236166
$ws.Range("O7").Value = 0.1068866440482915
$ws.Range("P7").Value = 0.1068866440482915
$ws.Range("Q7").Value = 394.2685491486638
$ws.Range("R7").Value = 3548.416942337974
$ws.Range("S7").Value = 0.05401491309797642
$ws.Range("T7").Value = 0.05401491309797642

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.894099000000001
$ws.Range("H8").Value = 26.682297
$ws.Range("I8").Value = 0.3732933164143983
$ws.Range("J8").Value = 0.3732933164143982
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 17.16653
$ws.Range("N8").Value = 51.49959
$ws.Range("O8").Value = 0.0560345397128279
$ws.Range("P8").Value = 0.0560345397128279
$ws.Range("Q8").Value = 152.68081730647
$ws.Range("R8").Value = 1374.12735575823
$ws.Range("S8").Value = 0.02091731916315583
$ws.Range("T8").Value = 0.02091731916315583

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.894099000000001
$ws.Range("H9").Value = 26.682297
$ws.Range("I9").Value = 0.3732933164143983
$ws.Range("J9").Value = 0.3732933164143982
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 256.4443053333333
$ws.Range("N9").Value = 769.332916
$ws.Range("O9").Value = 0.8370788162388805
$ws.Range("P9").Value = 0.8370788162388805
$ws.Range("Q9").Value = 2280.841039620895
$ws.Range("R9").Value = 20527.56935658805
$ws.Range("S9").Value = 0.3124759274140503
$ws.Range("T9").Value = 0.3124759274140503

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.894099000000001
$ws.Range("H10").Value = 26.682297
$ws.Range("I10").Value = 0.3732933164143983
$ws.Range("J10").Value = 0.3732933164143982
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 32.74538866666666
$ws.Range("N10").Value = 98.236166
$ws.Range("O10").Value = 0.1068866440482915
$ws.Range("P10").Value = 0.1068866440482915
$ws.Range("Q10").Value = 291.2407285948113
$ws.Range("R10").Value = 2621.166557353302
$ws.Range("S10").Value = 0.03990006983719205
$ws.Range("T10").Value = 0.03990006983719204

